$d = $word.ActiveDocument

# Keep literal straight quotes/apostrophes out of Word's "smart quotes" autocorrect
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$word.Options.AutoFormatReplaceQuotes = $false

# ---------------------------------------------------------------------------
# Change 1: insert a new "Meta description" paragraph right after the H1
# title paragraph ("Play De Magorum free and discover its magical features").
# ---------------------------------------------------------------------------

# Create a brand-new Normal-styled paragraph by splitting the (Normal-style)
# paragraph that currently sits right after the "Gameplay" heading, then
# relocate the "Gameplay" heading paragraph to sit after it, so the new blank
# paragraph lands directly under the H1 title without inheriting Heading1/
# Heading2 paragraph styles.
$pNormalAnchor = $d.Paragraphs.Item(3)
$rAnchor = $pNormalAnchor.Range
$rAnchor.Collapse(1)
$rAnchor.InsertParagraphBefore()

$pGameplay = $d.Paragraphs.Item(2)
$pGameplay.Range.Cut()

$pIfYoure = $d.Paragraphs.Item(3)
$pasteTarget = $d.Range($pIfYoure.Range.Start, $pIfYoure.Range.Start)
$pasteTarget.Paste()

# Paragraph 2 is now the new blank paragraph. Fill it with the meta
# description text, then bold just the "Meta description" label.
$metaPara = $d.Paragraphs.Item(2)
$metaRange = $metaPara.Range
$insertPos = $d.Range($metaRange.Start, $metaRange.Start)

$metaFullText = "Meta description: Read our review of De Magorum, play for free, and experience its magical symbols and Bonus Game feature with a high RTP and wide betting range."
$insertPos.InsertBefore($metaFullText)

$labelStart = $metaRange.Start
$labelEnd = $labelStart + ("Meta description").Length
$labelRange = $d.Range($labelStart, $labelEnd)
$labelRange.Font.Bold = 1

# ---------------------------------------------------------------------------
# Change 2: remove the duplicated bold "Play De Magorum free and discover
# its magical features" paragraph near the end of the document.
# ---------------------------------------------------------------------------

$oldBoldText = "Play De Magorum free and discover its magical features"
$dupPara = $null
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.Trim() -eq $oldBoldText -and $cand.Range.Font.Bold) {
        $dupPara = $cand
        break
    }
}
if ($dupPara -eq $null) {
    $dupPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
}
$dupPara.Range.Delete()

# ---------------------------------------------------------------------------
# Change 3: replace the text of the final (italic) paragraph with the new
# image-prompt copy, keeping its italic run formatting intact.
# ---------------------------------------------------------------------------

$lastPara = $d.Paragraphs.Last
$searchRange = $lastPara.Range.Duplicate
$oldItalicText = "Read our review of De Magorum, play for free, and experience its magical symbols and Bonus Game feature with a high RTP and wide betting range."
$newItalicText = 'Please create a cartoon-style image of a happy Maya warrior wearing glasses for Giocaonline''s slot game "De Magorum". The image should be engaging, vibrant, and playful, reflecting the excitement and energy of the game. Use warm colors and bold lines to create a friendly and approachable tone, and add magical elements such as books, potions, and crystal spheres to create a sense of mystery and fantasy. The happy Maya warrior should be in the center of the image, with a glowing staff held in his hand, and surrounded by the four magicians and their magical props. The Giocaonline logo should be prominently displayed at the bottom, and the text "De Magorum" should be written in an elegant and striking font.'

$found = $searchRange.Find.Execute($oldItalicText)
if ($found) {
    $searchRange.Text = $newItalicText
}

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
